$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2626439540334502
$ws.Range("C2").Value = 0.04467084256702947
$ws.Range("D2").Value = 0.02909498636159213
$ws.Range("E2").Value = 0.1616557201687172
$ws.Range("F2").Value = 0.6761016312668104
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.5142466153288936
$ws.Range("K2").Value = 0.2759508768186265
$ws.Range("M2").Value = 0.2237944772582168
$ws.Range("N2").Value = 1.393739124543001
$ws.Range("O2").Value = 2.30978727266762

$ws.Range("B3").Value = 0.2312670127796537
$ws.Range("C3").Value = 0.03948503063014641
$ws.Range("D3").Value = 0.02690771662449265
$ws.Range("E3").Value = 0.1506542875430057
$ws.Range("F3").Value = 0.6739632195842589
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.5177883982839724
$ws.Range("K3").Value = 0.2412588829804037
$ws.Range("M3").Value = 0.2014137900475674
$ws.Range("N3").Value = 1.409158980308042
$ws.Range("O3").Value = 2.315259363369634

$ws.Range("B4").Value = 0.2119921810154608
$ws.Range("C4").Value = 0.03628295256983449
$ws.Range("D4").Value = 0.02555262494102806
$ws.Range("E4").Value = 0.1440088162510591
$ws.Range("F4").Value = 0.6730606114817874
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.5202348931081247
$ws.Range("K4").Value = 0.2199201065871534
$ws.Range("M4").Value = 0.1877455117827509
$ws.Range("N4").Value = 1.41910945489443
$ws.Range("O4").Value = 2.319938753767772

$ws.Range("B5").Value = 0.204135634439325
$ws.Range("C5").Value = 0.03497359817400536
$ws.Range("D5").Value = 0.02499740016418173
$ws.Range("E5").Value = 0.1413280663308214
$ws.Range("F5").Value = 0.6727960059232245
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.5213002022678346
$ws.Range("K5").Value = 0.2112153133951438
$ws.Range("M5").Value = 0.1821940663906929
$ws.Range("N5").Value = 1.423285683620098
$ws.Range("O5").Value = 2.322177404874537

$ws.Range("B6").Value = 0.2028309590987192
$ws.Range("C6").Value = 0.03475591128278666
$ws.Range("D6").Value = 0.02490502431923147
$ws.Range("E6").Value = 0.1408845752026266
$ws.Range("F6").Value = 0.6727583028703208
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.52148122347959
$ws.Range("K6").Value = 0.20976935390172
$ws.Range("M6").Value = 0.1812733696292668
$ws.Range("N6").Value = 1.423986470400871
$ws.Range("O6").Value = 2.322569167673052

$ws.Range("B7").Value = 0.2118862319836694
$ws.Range("C7").Value = 0.03626531225560825
$ws.Range("D7").Value = 0.02554514915100725
$ws.Range("E7").Value = 0.1439725522974413
$ws.Range("F7").Value = 0.6730566249812284
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.5202489835371331
$ws.Range("K7").Value = 0.2198027468574963
$ws.Range("M7").Value = 0.1876705681960402
$ws.Range("N7").Value = 1.419165285756381
$ws.Range("O7").Value = 2.319967601763537

$ws.Range("B8").Value = 0.251827357537735
$ws.Range("C8").Value = 0.04288652456448006
$ws.Range("D8").Value = 0.02834334347544853
$ws.Range("E8").Value = 0.157839591933282
$ws.Range("F8").Value = 0.6752791407177909
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.5154113981853996
$ws.Range("K8").Value = 0.2639971570717137
$ws.Range("M8").Value = 0.2160623391304597
$ws.Range("N8").Value = 1.3989557066162
$ws.Range("O8").Value = 2.311400157977033

$ws.Range("B9").Value = 0.3300630182176576
$ws.Range("C9").Value = 0.0557271054068309
$ws.Range("D9").Value = 0.03373357094027796
$ws.Range("E9").Value = 0.1859112902867537
$ws.Range("F9").Value = 0.6828940875429765
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.5080824076853041
$ws.Range("K9").Value = 0.3503481696016308
$ws.Range("M9").Value = 0.272325380823446
$ws.Range("N9").Value = 1.363154859714989
$ws.Range("O9").Value = 2.305073390092559

$ws.Range("B10").Value = 0.3874736691849989
$ws.Range("C10").Value = 0.06507303656526631
$ws.Range("D10").Value = 0.03763359301412805
$ws.Range("E10").Value = 0.2070871077645933
$ws.Range("F10").Value = 0.6904764574293907
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.5040143433520115
$ws.Range("K10").Value = 0.4135852775733326
$ws.Range("M10").Value = 0.3140283333137504
$ws.Range("N10").Value = 1.339186956067675
$ws.Range("O10").Value = 2.306819930583089

$ws.Range("B11").Value = 0.4135734076582196
$ws.Range("C11").Value = 0.06930556936252685
$ws.Range("D11").Value = 0.03939454736510584
$ws.Range("E11").Value = 0.2168438565598834
$ws.Range("F11").Value = 0.6943580956012454
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.5024498549214016
$ws.Range("K11").Value = 0.4423064552839833
$ws.Range("M11").Value = 0.333081777984404
$ws.Range("N11").Value = 1.328790374165945
$ws.Range("O11").Value = 2.309005371899985

$ws.Range("B12").Value = 0.4234539353992943
$ws.Range("C12").Value = 0.0709055646759964
$ws.Range("D12").Value = 0.0400594538467729
$ws.Range("E12").Value = 0.2205565188960108
$ws.Range("F12").Value = 0.6958901599160399
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.5018985859326968
$ws.Range("K12").Value = 0.4531754927479881
$ws.Range("M12").Value = 0.3403087657652719
$ws.Range("N12").Value = 1.324926337391293
$ws.Range("O12").Value = 2.310033081221889

$ws.Range("B13").Value = 0.42132612258564
$ws.Range("C13").Value = 0.07056110108889868
$ws.Range("D13").Value = 0.03991634064090022
$ws.Range("E13").Value = 0.219756127254243
$ws.Range("F13").Value = 0.6955574373949815
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.5020154801200931
$ws.Range("K13").Value = 0.4508349721824345
$ws.Range("M13").Value = 0.338751777262658
$ws.Range("N13").Value = 1.32575528221219
$ws.Range("O13").Value = 2.309802842954412

$ws.Range("B14").Value = 0.4143863450374852
$ws.Range("C14").Value = 0.06943725776656606
$ws.Range("D14").Value = 0.03944928846368612
$ws.Range("E14").Value = 0.217148937788572
$ws.Range("F14").Value = 0.6944828934548326
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.5024036764404904
$ws.Range("K14").Value = 0.4432008015512565
$ws.Range("M14").Value = 0.3336761089549825
$ws.Range("N14").Value = 1.328471015558748
$ws.Range("O14").Value = 2.309085910598185

$ws.Range("B15").Value = 0.4101351419607226
$ws.Range("C15").Value = 0.0687485088501063
$ws.Range("D15").Value = 0.03916295337283771
$ws.Range("E15").Value = 0.2155543080411846
$ws.Range("F15").Value = 0.693832800797658
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.5026468201459799
$ws.Range("K15").Value = 0.4385237194111937
$ws.Range("M15").Value = 0.330568659983399
$ws.Range("N15").Value = 1.330143979919693
$ws.Range("O15").Value = 2.3086728351023

$ws.Range("B16").Value = 0.3857676135763484
$ws.Range("C16").Value = 0.06479604494694513
$ws.Range("D16").Value = 0.037518242391954
$ws.Range("E16").Value = 0.2064519885075526
$ws.Range("F16").Value = 0.6902314835554364
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.5041223457639923
$ws.Range("K16").Value = 0.4117073206468547
$ws.Range("M16").Value = 0.3127848063828509
$ws.Range("N16").Value = 1.339876597074315
$ws.Range("O16").Value = 2.306705105875835

$ws.Range("B17").Value = 0.3708143157547568
$ws.Range("C17").Value = 0.06236644434325456
$ws.Range("D17").Value = 0.03650586517537846
$ws.Range("E17").Value = 0.2008998488276035
$ws.Range("F17").Value = 0.6881329327389807
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.5051008290861496
$ws.Range("K17").Value = 0.3952442840498804
$ws.Range("M17").Value = 0.3018961272557021
$ws.Range("N17").Value = 1.345977037666801
$ws.Range("O17").Value = 2.305854301370744

$ws.Range("B18").Value = 0.3622120355556433
$ws.Range("C18").Value = 0.06096721865567645
$ws.Range("D18").Value = 0.03592233347639251
$ws.Range("E18").Value = 0.1977180424039346
$ws.Range("F18").Value = 0.6869665997379855
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.5056905536538139
$ws.Range("K18").Value = 0.3857709079174185
$ws.Range("M18").Value = 0.2956410188877854
$ws.Range("N18").Value = 1.349533522306449
$ws.Range("O18").Value = 2.305495855276575

$ws.Range("B19").Value = 0.3592992016785104
$ws.Range("C19").Value = 0.06049315986844306
$ws.Range("D19").Value = 0.035724547847515
$ws.Range("E19").Value = 0.1966427301779561
$ws.Range("F19").Value = 0.6865786894271153
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.5058948480937175
$ws.Range("K19").Value = 0.3825626662031709
$ws.Range("M19").Value = 0.2935244818015477
$ws.Range("N19").Value = 1.350745872198591
$ws.Range("O19").Value = 2.305396972678892

$ws.Range("B20").Value = 0.3724062825837677
$ws.Range("C20").Value = 0.06262526439817862
$ws.Range("D20").Value = 0.03661376293786134
$ws.Range("E20").Value = 0.2014896777180439
$ws.Range("F20").Value = 0.6883521149958085
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.5049938808297583
$ws.Range("K20").Value = 0.3969972480895194
$ws.Range("M20").Value = 0.3030544414522538
$ws.Range("N20").Value = 1.345322700819675
$ws.Range("O20").Value = 2.305931320993409

$ws.Range("B21").Value = 0.4164248078425601
$ws.Range("C21").Value = 0.06976743340248959
$ws.Range("D21").Value = 0.03958652561502163
$ws.Range("E21").Value = 0.2179142426377325
$ws.Range("F21").Value = 0.6947968257549064
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.5022885362317311
$ws.Range("K21").Value = 0.4454433376427858
$ws.Range("M21").Value = 0.335166633937007
$ws.Range("N21").Value = 1.327671359052218
$ws.Range("O21").Value = 2.309291058616964

$ws.Range("B22").Value = 0.4451764668426677
$ws.Range("C22").Value = 0.07441907010725402
$ws.Range("D22").Value = 0.04151814064772452
$ws.Range("E22").Value = 0.2287536149551812
$ws.Range("F22").Value = 0.6993712015774207
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.5007604139089992
$ws.Range("K22").Value = 0.4770643476340126
$ws.Range("M22").Value = 0.3562229420233081
$ws.Range("N22").Value = 1.316560279304946
$ws.Range("O22").Value = 2.312653406006547

$ws.Range("B23").Value = 0.4298328793688313
$ws.Range("C23").Value = 0.07193789940063766
$ws.Range("D23").Value = 0.04048824220452474
$ws.Range("E23").Value = 0.2229587684464107
$ws.Range("F23").Value = 0.6968966142794315
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.5015540341581257
$ws.Range("K23").Value = 0.4601915572914095
$ws.Range("M23").Value = 0.3449784721319489
$ws.Range("N23").Value = 1.322451547424496
$ws.Range("O23").Value = 2.31075207914219

$ws.Range("B24").Value = 0.3716865712581807
$ws.Range("C24").Value = 0.0625082593745816
$ws.Range("D24").Value = 0.03656498703747957
$ws.Range("E24").Value = 0.2012229843448026
$ws.Range("F24").Value = 0.6882528976267963
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.505042147469382
$ws.Range("K24").Value = 0.3962047596664888
$ws.Range("M24").Value = 0.3025307522238023
$ws.Range("N24").Value = 1.345618373069105
$ws.Range("O24").Value = 2.305896093327817

$ws.Range("B25").Value = 0.3089092496148567
$ws.Range("C25").Value = 0.05226881286985474
$ws.Range("D25").Value = 0.03228586489769469
$ws.Range("E25").Value = 0.1782214503703514
$ws.Range("F25").Value = 0.6804851550020885
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.5098339870783803
$ws.Range("K25").Value = 0.3270229531894415
$ws.Range("M25").Value = 0.2570410035936774
$ws.Range("N25").Value = 1.372429904822191
$ws.Range("O25").Value = 2.305662644113539

